# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the "Status" column moves
# from "Handed back: in sync with en-US" to "Ready for handoff", and the
# associated timestamp columns are refreshed. The zh-cn/de-de column that
# used to show the "Status" text also gets auto-sized narrower to fit the
# new (shorter) value.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed "Latest HO Xliff Generate Date" / handoff timestamps
$wsOverview.Range("G2").Value = "2016-08-19 02:56:16"
$wsDeDe.Range("H2").Value     = "2016-08-19 02:56:16"
$wsZhCn.Range("H2").Value     = "2016-08-19 02:56:11"

# --- Column width: the Status column narrows now that the text is shorter
# ("Ready for handoff" vs "Handed back: in sync with en-US"), re-sized to
# roughly 17.22 characters wide (from ~29.98). ColumnWidth snaps to the
# nearest whole-pixel increment, so feed it a value inside that pixel
# bucket to land as close as possible to the target width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33
